$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("OC", 4.02, 4.97),
    @("MC", 5.43, 5.69),
    @("SPC", 6.58, 23.4),
    @("8319", 2.99, 3.06),
    @("Erin", 0, 1.1),
    @("143", 1.95, 4.95),
    @("11495", 10.53, 30.77),
    @("1371", 3.01, 5.21)
)

$row = 2
foreach ($r in $data) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r[0]
    $cellA.ClearFormats()
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
